# Weekly update: insert a new daily price record for
# "Pepino ensalada" (Vega Central Mapocho de Santiago) above the existing
# row 387. This pushes all the following rows down by one (387 -> 436)
# and grows the used range from A1:R435 to A1:R436.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 387; rows 387..435 shift down to 388..436.
$ws.Rows("387:387").Insert()

# Populate the freshly inserted row with the new record.
$ws.Range("A387").Value = 9
$ws.Range("B387").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C387").Value = "Metropolitana"
$ws.Range("D387").Value2 = 45124
$ws.Range("E387").Value = 13
$ws.Range("F387").Value = 100112043
$ws.Range("G387").Value = "Pepino ensalada"
$ws.Range("H387").Value = "Sin especificar"
$ws.Range("I387").Value = "Primera"
$ws.Range("J387").Value = 70
$ws.Range("K387").Value = 12000
$ws.Range("L387").Value = 13000
$ws.Range("M387").Value = 12500
$ws.Range("N387").Value = "$/caja 60 unidades"
$ws.Range("O387").Value = "Región de Arica y Parinacota"
$ws.Range("P387").Value = 208
$ws.Range("Q387").Value = 60
$ws.Range("R387").Value = "Hortaliza"
